$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 13:20"

# Row 6: numbers updated
$ws.Range("B6").Value = 3315017
$ws.Range("C6").Value = 7268
$ws.Range("D6").Value = 2524539
$ws.Range("E6").Value = 729826

# Row 15: numbers updated
$ws.Range("B15").Value = 367796
$ws.Range("C15").Value = 2190
$ws.Range("D15").Value = 316638
$ws.Range("E15").Value = 30021
$ws.Range("G15").Value = 117
$ws.Range("H15").Value = 21137

# Row 40: country -> Rumania, numbers updated
$ws.Range("A40").Value = "Rumania"
$ws.Range("B40").Value = 83150
$ws.Range("C40").Value = 1504
$ws.Range("D40").Value = 36677
$ws.Range("E40").Value = 43014
$ws.Range("G40").Value = 38
$ws.Range("H40").Value = 3459

# Row 41: country -> Belgica, numbers updated
$ws.Range("A41").Value = "Belgica"
$ws.Range("B41").Value = 83030
$ws.Range("C41").Value = 583
$ws.Range("D41").Value = 18331
$ws.Range("E41").Value = 54820
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 9879

# Row 42: country -> Kuwait, numbers updated
$ws.Range("A42").Value = "Kuwait"
$ws.Range("B42").Value = 82271
$ws.Range("D42").Value = 73906
$ws.Range("E42").Value = 7844
$ws.Range("H42").Value = 521

# Row 61: numbers updated
$ws.Range("B61").Value = 41006
$ws.Range("C61").Value = 361
$ws.Range("E61").Value = 4203

# Row 66: country -> Nepal, numbers updated
$ws.Range("A66").Value = "Nepal"
$ws.Range("B66").Value = 35529
$ws.Range("C66").Value = 1111
$ws.Range("D66").Value = 20073
$ws.Range("E66").Value = 15273
$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 183

# Row 67: country -> Moldavia, numbers updated
$ws.Range("A67").Value = "Moldavia"
$ws.Range("B67").Value = 34982
$ws.Range("D67").Value = 24156
$ws.Range("E67").Value = 9859
$ws.Range("H67").Value = 967

# Row 72: numbers updated
$ws.Range("D72").Value = 20366
$ws.Range("E72").Value = 4384

# Row 75: numbers updated
$ws.Range("B75").Value = 20677
$ws.Range("C75").Value = 522
$ws.Range("D75").Value = 14194
$ws.Range("E75").Value = 6342
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 141

# Row 129: numbers updated
$ws.Range("B129").Value = 2755
$ws.Range("C129").Value = 33
$ws.Range("D129").Value = 2199
$ws.Range("E129").Value = 423

# Row 130: country -> Gambia, numbers updated
$ws.Range("A130").Value = "Gambia"
$ws.Range("B130").Value = 2743
$ws.Range("C130").Value = 35
$ws.Range("D130").Value = 638
$ws.Range("E130").Value = 2012
$ws.Range("H130").Value = 93

# Row 131: country -> Mali, numbers updated
$ws.Range("A131").Value = "Mali"
$ws.Range("B131").Value = 2717
$ws.Range("D131").Value = 2052
$ws.Range("E131").Value = 539
$ws.Range("H131").Value = 126

# Row 144: country -> Malta, numbers updated
$ws.Range("A144").Value = "Malta"
$ws.Range("B144").Value = 1788
$ws.Range("C144").Value = 37
$ws.Range("D144").Value = 1121
$ws.Range("E144").Value = 657
$ws.Range("H144").Value = 10

# Row 145: country -> Aruba, numbers updated
$ws.Range("A145").Value = "Aruba"
$ws.Range("B145").Value = 1760
$ws.Range("D145").Value = 587
$ws.Range("E145").Value = 1165
$ws.Range("H145").Value = 8

# Row 146: country -> Jordania, numbers updated
$ws.Range("A146").Value = "Jordania"
$ws.Range("B146").Value = 1756
$ws.Range("D146").Value = 1355
$ws.Range("E146").Value = 386
$ws.Range("H146").Value = 15

# Row 162: numbers updated
$ws.Range("B162").Value = 1036
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 637
$ws.Range("E162").Value = 369

# Row 202: country -> Timor Oriental, numbers updated
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("B202").Value = 27
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 25
$ws.Range("E202").Value = 2
$ws.Range("H202").Value = 0

# Row 203: country -> Islas Virgenes Britanicas, numbers updated
$ws.Range("A203").Value = "Islas Virgenes Britanicas"
$ws.Range("D203").Value = 8
$ws.Range("E203").Value = 17
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 1

# Row 204: country -> Santa Lucia
$ws.Range("A204").Value = "Santa Lucia"
